$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.064450892247345
$ws.Range("D2").Value = 1.065609612849498
$ws.Range("E2").Value = 1.067193751013686
$ws.Range("F2").Value = 1.077325582023301
$ws.Range("I2").Value = 1.042971575057575
$ws.Range("J2").Value = 1.069411067594375
$ws.Range("K2").Value = 1.068322391128056
$ws.Range("L2").Value = 1.069902269760689
$ws.Range("M2").Value = 1.08000717781526
$ws.Range("N2").Value = 1.070929753973523
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.066195067267269
$ws.Range("D3").Value = 1.066972285970766
$ws.Range("E3").Value = 1.068606784435461
$ws.Range("F3").Value = 1.078850171154463
$ws.Range("I3").Value = 1.04336489812504
$ws.Range("J3").Value = 1.070806695590494
$ws.Range("K3").Value = 1.069499129451866
$ws.Range("L3").Value = 1.071129554312495
$ws.Range("M3").Value = 1.081347714112699
$ws.Range("N3").Value = 1.072327363921477
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.067321291634983
$ws.Range("D4").Value = 1.067851767729811
$ws.Range("E4").Value = 1.069519809211658
$ws.Range("F4").Value = 1.079835006252076
$ws.Range("I4").Value = 1.043617024183455
$ws.Range("J4").Value = 1.071707062227843
$ws.Range("K4").Value = 1.070257764311699
$ws.Range("L4").Value = 1.071921860106372
$ws.Range("M4").Value = 1.082212947302763
$ws.Range("N4").Value = 1.073229009182725
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.067794200988787
$ws.Range("D5").Value = 1.068220970019979
$ws.Range("E5").Value = 1.069903341186758
$ws.Range("F5").Value = 1.080248639135252
$ws.Range("I5").Value = 1.043722451393512
$ws.Range("J5").Value = 1.072084940515972
$ws.Range("K5").Value = 1.070576033987635
$ws.Range("L5").Value = 1.072254515096612
$ws.Range("M5").Value = 1.082576177382631
$ws.Range("N5").Value = 1.073607424101366
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.067873572354592
$ws.Range("D6").Value = 1.068282929779945
$ws.Range("E6").Value = 1.069967720322726
$ws.Range("F6").Value = 1.080318067214764
$ws.Range("I6").Value = 1.043740119955221
$ws.Range("J6").Value = 1.072148350947791
$ws.Range("K6").Value = 1.070629434380879
$ws.Range("L6").Value = 1.07231034431584
$ws.Range("M6").Value = 1.08263713542625
$ws.Range("N6").Value = 1.073670924583271
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.067327612837435
$ws.Range("D7").Value = 1.067856703106252
$ws.Range("E7").Value = 1.069524935166446
$ws.Range("F7").Value = 1.079840534761261
$ws.Range("I7").Value = 1.043618435128738
$ws.Range("J7").Value = 1.071712113940551
$ws.Range("K7").Value = 1.070262019630093
$ws.Range("L7").Value = 1.071926306739113
$ws.Range("M7").Value = 1.082217802808187
$ws.Range("N7").Value = 1.073234068069445
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.065040843572682
$ws.Range("D8").Value = 1.066070607724848
$ws.Range("E8").Value = 1.067671566849686
$ws.Range("F8").Value = 1.077841176797838
$ws.Range("I8").Value = 1.043104995248059
$ws.Range("J8").Value = 1.069883291326747
$ws.Range("K8").Value = 1.068720658820417
$ws.Range("L8").Value = 1.070317419692332
$ws.Range("M8").Value = 1.08046067572365
$ws.Range("N8").Value = 1.07140264831776
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.060992507369162
$ws.Range("D9").Value = 1.062905576518285
$ws.Range("E9").Value = 1.064395354024431
$ws.Range("F9").Value = 1.074304790687384
$ws.Range("I9").Value = 1.042181871641638
$ws.Range("J9").Value = 1.066639555643625
$ws.Range("K9").Value = 1.065982812553314
$ws.Range("L9").Value = 1.067468002295828
$ws.Range("M9").Value = 1.07734727592436
$ws.Range("N9").Value = 1.068154306157905
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.058280193879214
$ws.Range("D10").Value = 1.060783086896929
$ws.Range("E10").Value = 1.062203743600619
$ws.Range("F10").Value = 1.071937685017994
$ws.Range("I10").Value = 1.041553893112529
$ws.Range("J10").Value = 1.064462230889153
$ws.Range("K10").Value = 1.0641424141411
$ws.Range("L10").Value = 1.065558259577978
$ws.Range("M10").Value = 1.075259612064182
$ws.Range("N10").Value = 1.06597388935254
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.057102367009304
$ws.Range("D11").Value = 1.059860941042275
$ws.Range("E11").Value = 1.061252865000406
$ws.Range("F11").Value = 1.070910311246571
$ws.Range("I11").Value = 1.041278945126352
$ws.Range("J11").Value = 1.063515764156348
$ws.Range("K11").Value = 1.063341785026334
$ws.Range("L11").Value = 1.064728809053839
$ws.Range("M11").Value = 1.074352644697872
$ws.Range("N11").Value = 1.065026078528413
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.056664346553844
$ws.Range("D12").Value = 1.059517939158198
$ws.Range("E12").Value = 1.060899371322379
$ws.Range("F12").Value = 1.070528326458195
$ws.Range("I12").Value = 1.041176358021405
$ws.Range("J12").Value = 1.063163640400915
$ws.Range("K12").Value = 1.063043825836257
$ws.Range("L12").Value = 1.064420326315111
$ws.Range("M12").Value = 1.074015295832093
$ws.Range("N12").Value = 1.064673454716858
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.056758327354516
$ws.Range("D13").Value = 1.059591535938607
$ws.Range("E13").Value = 1.060975210429396
$ws.Range("F13").Value = 1.070610280533237
$ws.Range("I13").Value = 1.041198384148644
$ws.Range("J13").Value = 1.063239197858178
$ws.Range("K13").Value = 1.063107765049659
$ws.Range("L13").Value = 1.064486514613916
$ws.Range("M13").Value = 1.074087679336138
$ws.Range("N13").Value = 1.064749119474378
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.057066170864684
$ws.Range("D14").Value = 1.059832598171117
$ws.Range("E14").Value = 1.061223651179703
$ws.Range("F14").Value = 1.070878743929957
$ws.Range("I14").Value = 1.041270474635207
$ws.Range("J14").Value = 1.063486669071723
$ws.Range("K14").Value = 1.063317167312445
$ws.Range("L14").Value = 1.064703317741736
$ws.Range("M14").Value = 1.074324768798612
$ws.Range("N14").Value = 1.064996942125431
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.057255773620772
$ws.Range("D15").Value = 1.059981061180941
$ws.Range("E15").Value = 1.061376684329006
$ws.Range("F15").Value = 1.071044103388688
$ws.Range("I15").Value = 1.041314830994321
$ws.Range("J15").Value = 1.063639069163215
$ws.Range("K15").Value = 1.063446111130969
$ws.Range("L15").Value = 1.064836845623619
$ws.Range("M15").Value = 1.074470786054958
$ws.Range("N15").Value = 1.065149558642533
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.058358289876718
$ws.Range("D16").Value = 1.060844220556871
$ws.Range("E16").Value = 1.062266809462144
$ws.Range("F16").Value = 1.072005816845873
$ws.Range("I16").Value = 1.041572076336547
$ws.Range("J16").Value = 1.064524966376775
$ws.Range("K16").Value = 1.064195469909237
$ws.Range("L16").Value = 1.065613253599584
$ws.Range("M16").Value = 1.075319740476339
$ws.Range("N16").Value = 1.066036713931751
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.05904895396504
$ws.Range("D17").Value = 1.061384821561388
$ws.Range("E17").Value = 1.062824646543345
$ws.Range("F17").Value = 1.072608423184706
$ws.Range("I17").Value = 1.041732625685442
$ws.Range("J17").Value = 1.065079674829486
$ws.Range("K17").Value = 1.064664518470086
$ws.Range("L17").Value = 1.066099592721625
$ws.Range("M17").Value = 1.07585145766996
$ws.Range("N17").Value = 1.06659221013409
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.059451481858995
$ws.Range("D18").Value = 1.061699847222154
$ws.Range("E18").Value = 1.063149841143548
$ws.Range("F18").Value = 1.072959682359865
$ws.Range("I18").Value = 1.04182597941495
$ws.Range("J18").Value = 1.065402873554296
$ws.Range("K18").Value = 1.064937747907478
$ws.Range("L18").Value = 1.066383023664041
$ws.Range("M18").Value = 1.076161311062353
$ws.Range("N18").Value = 1.066915867838161
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.059588678852857
$ws.Range("D19").Value = 1.061807212658503
$ws.Range("E19").Value = 1.063260693474224
$ws.Range("F19").Value = 1.073079413845694
$ws.Range("I19").Value = 1.041857761220637
$ws.Range("J19").Value = 1.065513016408773
$ws.Range("K19").Value = 1.065030851561845
$ws.Range("L19").Value = 1.066479625451228
$ws.Range("M19").Value = 1.076266914494938
$ws.Range("N19").Value = 1.067026167108124
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.05897488595224
$ws.Range("D20").Value = 1.061326851068207
$ws.Range("E20").Value = 1.062764814835934
$ws.Range("F20").Value = 1.072543793195529
$ws.Range("I20").Value = 1.041715430491112
$ws.Range("J20").Value = 1.065020196498351
$ws.Range("K20").Value = 1.064614231170318
$ws.Range("L20").Value = 1.066047438280973
$ws.Range("M20").Value = 1.075794439340529
$ws.Range("N20").Value = 1.066532647336903
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.056975533191273
$ws.Range("D21").Value = 1.05976162457713
$ws.Range("E21").Value = 1.061150499796414
$ws.Range("F21").Value = 1.070799698525772
$ws.Range("I21").Value = 1.041249258492403
$ws.Range("J21").Value = 1.063413810633013
$ws.Range("K21").Value = 1.063255519380495
$ws.Range("L21").Value = 1.064639485407562
$ws.Range("M21").Value = 1.074254964669066
$ws.Range("N21").Value = 1.06492398021938
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.055715426282019
$ws.Range("D22").Value = 1.058774745855024
$ws.Range("E22").Value = 1.06013380386431
$ws.Range("F22").Value = 1.069700956389086
$ws.Range("I22").Value = 1.040953498919124
$ws.Range("J22").Value = 1.062400543482374
$ws.Range("K22").Value = 1.06239794095418
$ws.Range("L22").Value = 1.063752000666962
$ws.Range("M22").Value = 1.073284364337086
$ws.Range("N22").Value = 1.063909274113165
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.05638372571939
$ws.Range("D23").Value = 1.059298174146649
$ws.Range("E23").Value = 1.060672939273791
$ws.Range("F23").Value = 1.070283629101694
$ws.Range("I23").Value = 1.041110540058088
$ws.Range("J23").Value = 1.062938009562259
$ws.Range("K23").Value = 1.062852876009774
$ws.Range("L23").Value = 1.064222689445336
$ws.Range("M23").Value = 1.073799154805199
$ws.Range("N23").Value = 1.064447503456531
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.059008355102472
$ws.Range("D24").Value = 1.061353046358861
$ws.Range("E24").Value = 1.062791850772582
$ws.Range("F24").Value = 1.072572997419222
$ws.Range("I24").Value = 1.041723201161519
$ws.Range("J24").Value = 1.065047073287721
$ws.Range("K24").Value = 1.064636954944716
$ws.Range("L24").Value = 1.066071005376934
$ws.Range("M24").Value = 1.075820204359991
$ws.Range("N24").Value = 1.066559562294397
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.062041405181099
$ws.Range("D25").Value = 1.06372596880147
$ws.Range("E25").Value = 1.065243611834851
$ws.Range("F25").Value = 1.075220664118294
$ws.Range("I25").Value = 1.042422720175164
$ws.Range("J25").Value = 1.067480706663572
$ws.Range("K25").Value = 1.066693245680634
$ws.Range("L25").Value = 1.068206395119626
$ws.Range("M25").Value = 1.07815424972695
$ws.Range("N25").Value = 1.068996651708782
